$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "NIDHI"
$ws.Range("B7").Value = "Math"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 2

$ws.Range("A8").Value = "J"
$ws.Range("B8").Value = "Python"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2
